$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2 (shifts existing rows 2-18 down to 3-19)
$ws.Rows("2:2").Insert()

# Clear any auto-copied formatting on the new row, then apply the same style
# pattern used by the other data rows (date style on column A only).
$ws.Range("A2:E2").ClearFormats()
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Write the refreshed data values (new row plus recalculated rows below it)
$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 2007
$ws.Range("C2").Value = 7.226520411029047
$ws.Range("D2").Value = 2008
$ws.Range("E2").Value = 13.12477391005418

$ws.Range("A3").Value = 39765
$ws.Range("B3").Value = 2008
$ws.Range("C3").Value = 4.268860212333636
$ws.Range("D3").Value = 2009
$ws.Range("E3").Value = 13.08276537368067

$ws.Range("A4").Value = 40130
$ws.Range("B4").Value = 2009
$ws.Range("C4").Value = -7.266312015249776
$ws.Range("D4").Value = 2010
$ws.Range("E4").Value = 12.31225042954258

$ws.Range("A5").Value = 40494
$ws.Range("B5").Value = 2010
$ws.Range("C5").Value = 6.958243460951929
$ws.Range("D5").Value = 2011
$ws.Range("E5").Value = 12.04357532583245

$ws.Range("A6").Value = 40862
$ws.Range("B6").Value = 2011
$ws.Range("C6").Value = 9.469137444079934
$ws.Range("D6").Value = 2012
$ws.Range("E6").Value = 10.6992064911972

$ws.Range("A7").Value = 41228
$ws.Range("B7").Value = 2012
$ws.Range("C7").Value = 3.358206407534947
$ws.Range("D7").Value = 2013
$ws.Range("E7").Value = 4.390489499870109

$ws.Range("A8").Value = 41592
$ws.Range("B8").Value = 2013
$ws.Range("C8").Value = 0.3081076735359067
$ws.Range("D8").Value = 2014
$ws.Range("E8").Value = 3.50243535103556

$ws.Range("A9").Value = 41957
$ws.Range("B9").Value = 2014
$ws.Range("C9").Value = 3.901355411819707
$ws.Range("D9").Value = 2015
$ws.Range("E9").Value = 6.143002545701282

$ws.Range("A10").Value = 42321
$ws.Range("B10").Value = 2015
$ws.Range("C10").Value = 5.331683351557981
$ws.Range("D10").Value = 2016
$ws.Range("E10").Value = 4.555278923792572

$ws.Range("A11").Value = 42689
$ws.Range("B11").Value = 2016
$ws.Range("C11").Value = 3.254758369308375
$ws.Range("D11").Value = 2017
$ws.Range("E11").Value = 0.9515943257393689

$ws.Range("A12").Value = 43053
$ws.Range("B12").Value = 2017
$ws.Range("C12").Value = 5.246209615995667
$ws.Range("D12").Value = 2018
$ws.Range("E12").Value = 4.251116704684899

$ws.Range("A13").Value = 43418
$ws.Range("B13").Value = 2018
$ws.Range("C13").Value = 4.86255966374296
$ws.Range("D13").Value = 2019
$ws.Range("E13").Value = 4.992093705734701

$ws.Range("A14").Value = 43783
$ws.Range("B14").Value = 2019
$ws.Range("C14").Value = 2.764740011159428
$ws.Range("D14").Value = 2020
$ws.Range("E14").Value = 0.3611963426346065

$ws.Range("A15").Value = 44159
$ws.Range("B15").Value = 2020
$ws.Range("C15").Value = -7.260793671746435
$ws.Range("D15").Value = 2021
$ws.Range("E15").Value = 2.387971016884638

$ws.Range("A16").Value = 44525
$ws.Range("B16").Value = 2021
$ws.Range("C16").Value = 4.097586525396268
$ws.Range("D16").Value = 2022
$ws.Range("E16").Value = 3.347989317130651

$ws.Range("A17").Value = 44890
$ws.Range("B17").Value = 2022
$ws.Range("C17").Value = 7.824284864703746
$ws.Range("D17").Value = 2023
$ws.Range("E17").Value = -2.122471977790918

$ws.Range("A18").Value = 45254
$ws.Range("B18").Value = 2023
$ws.Range("C18").Value = -1.24502235313334
$ws.Range("D18").Value = 2024
$ws.Range("E18").Value = -4.678511595261359

$ws.Range("A19").Value = 45618
$ws.Range("B19").Value = 2024
$ws.Range("C19").Value = -1.735114423676209
$ws.Range("D19").Value = 2025
$ws.Range("E19").Value = 1.758477003221981

Write-Host "edit complete"
